# ---------------------------------------------------------------------------
# feat: add 2022-Q1 data
#
# Before:  Sheets = [ "2021-Q4", "总计" ]
# After:   Sheets = [ "2021-Q4", "2022-Q1", "总计" ]
#   - "2021-Q4"  unchanged (per-fund holdings table)
#   - "2022-Q1"  NEW per-fund holdings table (same shape as "2021-Q4")
#   - "总计"     rebuilt summary table: a new top row for 2022-Q1 is added,
#                the existing 2021-Q4 summary row shifts down one row.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# xlPasteFormats -4122 ; used to copy over cell formatting without values.
$xlPasteFormats = -4122

# Helper: force a PowerShell string to be stored as literal text (not
# auto-coerced to a number), while leaving the cell's existing style/number
# format untouched (writing the value normally, then setting NumberFormat=@
# would stamp a brand-new style index onto the cell).
function Set-TextValue($ws, $addr, $text, $blankRef) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $text
    $blankRef.Copy()
    $cell.PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------------
# Step 1: rename the existing "总计" sheet (sheetId 2) to "2022-Q1", and make
# an in-workbook duplicate of it (placed immediately after) that becomes the
# new "总计" sheet (sheetId 3) -- this preserves the original 2-sheet file's
# sheetId/rId numbering exactly the way Excel would when inserting a sheet.
# ---------------------------------------------------------------------------
$wsOldTotal = $wb.Worksheets.Item("总计")
$wsOldTotal.Copy()
$wsNewTotal = $wb.Worksheets.Item(3)

$wsOldTotal.Name = "2022-Q1"
$wsNewTotal.Name = "总计"

$wsQ4 = $wb.Worksheets.Item("2021-Q4")
$wsQ1 = $wb.Worksheets.Item("2022-Q1")
$wsTotal = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------------
# Step 2: reshape "2022-Q1" from the old 4-column "总计" layout into the
# 8-column per-fund holdings layout, by copying formatting from "2021-Q4"
# (header row B1:H1, and the bold/bordered index column A2:A5).
# ---------------------------------------------------------------------------
$wsQ4.Range("B1:H1").Copy()
$wsQ1.Range("B1:H1").PasteSpecial($xlPasteFormats)
$wsQ4.Range("A2:A5").Copy()
$wsQ1.Range("A2:A5").PasteSpecial($xlPasteFormats)

$blankQ1 = $wsQ1.Range("Z100")

# Header row
$wsQ1.Range("B1").Value = "基金代码"
$wsQ1.Range("C1").Value = "基金名称"
$wsQ1.Range("D1").Value = "基金规模"
$wsQ1.Range("E1").Value = "股票总仓位"
$wsQ1.Range("F1").Value = "仓位占比"
$wsQ1.Range("G1").Value = "持有市值(亿元)"
$wsQ1.Range("H1").Value = "仓位排名"

# Row 2: 010936 / 交银施罗德均衡成长一年持有期混合A
$wsQ1.Range("A2").Value = 0
Set-TextValue $wsQ1 "B2" "010936" $blankQ1
Set-TextValue $wsQ1 "C2" "交银施罗德均衡成长一年持有期混合A" $blankQ1
Set-TextValue $wsQ1 "D2" "96.27" $blankQ1
Set-TextValue $wsQ1 "E2" "87.48" $blankQ1
Set-TextValue $wsQ1 "F2" "2.94" $blankQ1
Set-TextValue $wsQ1 "G2" "2.8303" $blankQ1
$wsQ1.Range("H2").Value = 10

# Row 3: 519704 / 交银先进制造混合
$wsQ1.Range("A3").Value = 1
Set-TextValue $wsQ1 "B3" "519704" $blankQ1
Set-TextValue $wsQ1 "C3" "交银先进制造混合" $blankQ1
Set-TextValue $wsQ1 "D3" "75.33" $blankQ1
Set-TextValue $wsQ1 "E3" "86.31" $blankQ1
Set-TextValue $wsQ1 "F3" "2.84" $blankQ1
Set-TextValue $wsQ1 "G3" "2.1394" $blankQ1
$wsQ1.Range("H3").Value = 10

# Row 4: 009402 / 交银施罗德启明混合
$wsQ1.Range("A4").Value = 2
Set-TextValue $wsQ1 "B4" "009402" $blankQ1
Set-TextValue $wsQ1 "C4" "交银施罗德启明混合" $blankQ1
Set-TextValue $wsQ1 "D4" "51.54" $blankQ1
Set-TextValue $wsQ1 "E4" "82.62" $blankQ1
Set-TextValue $wsQ1 "F4" "2.59" $blankQ1
Set-TextValue $wsQ1 "G4" "1.3349" $blankQ1
$wsQ1.Range("H4").Value = 10

# Row 5: 010937 / 交银施罗德均衡成长一年持有期混合C
$wsQ1.Range("A5").Value = 3
Set-TextValue $wsQ1 "B5" "010937" $blankQ1
Set-TextValue $wsQ1 "C5" "交银施罗德均衡成长一年持有期混合C" $blankQ1
Set-TextValue $wsQ1 "D5" "2.50" $blankQ1
Set-TextValue $wsQ1 "E5" "87.48" $blankQ1
Set-TextValue $wsQ1 "F5" "2.94" $blankQ1
Set-TextValue $wsQ1 "G5" "0.0735" $blankQ1
$wsQ1.Range("H5").Value = 10

$blankQ1.ClearContents()

# ---------------------------------------------------------------------------
# Step 3: rebuild "总计" - push the existing 2021-Q4 summary row to row 3
# (index 1) and write the new 2022-Q1 summary into row 2 (index 0).
# ---------------------------------------------------------------------------
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial($xlPasteFormats)

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2021-Q4"
$wsTotal.Range("C3").Value = 4
$wsTotal.Range("D3").Value = 10.52

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 4
$wsTotal.Range("D2").Value = 6.38
